# Update "想去人数" (F column) values on the 展览 (sheet1) and 全部类型 (sheet4) sheets
# to reflect the newly generated output numbers.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 216
$wsExhibition.Range("F10").Value = 45
$wsExhibition.Range("F11").Value = 7013
$wsExhibition.Range("F13").Value = 401
$wsExhibition.Range("F14").Value = 3454
$wsExhibition.Range("F15").Value = 255
$wsExhibition.Range("F16").Value = 448
$wsExhibition.Range("F18").Value = 582
$wsExhibition.Range("F19").Value = 60

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 216
$wsAll.Range("F12").Value = 45
$wsAll.Range("F14").Value = 7013
$wsAll.Range("F17").Value = 401
$wsAll.Range("F18").Value = 3454
$wsAll.Range("F19").Value = 255
$wsAll.Range("F20").Value = 448
$wsAll.Range("F22").Value = 582
$wsAll.Range("F23").Value = 60
